# Minor fixes after Kams review
# Applies the OOXML-diff described changes through the Word object model:
#  - merge runs that were previously split around proofErr spell-check
#    markers (no visible text change, just cleanup of the XML)
#  - several small wording/content tweaks
#  - turn the "CloudFix" mention in the summary bullet into a hyperlink
#  - tidy up the footer sentence into a single run

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# --- TECHNICAL EXPERTISE section: merge runs split by spell-check markers ---
Replace-Text "Java, Groovy, Javascript, Python" "Java, Groovy, Javascript, Python"
Replace-Text "API Gateway, AppSync, Lambda, EventBridge, Step Functions" "API Gateway, AppSync, Lambda, EventBridge, Step Functions"
Replace-Text "REST, GraphQL" "REST, GraphQL"
Replace-Text "Lucidchart, Draw.io" "Lucidchart, Draw.io"
Replace-Text "Quicksight, Grafana" "Quicksight, Grafana"

# --- CERTIFICATES section: merge the freeCodeCamp hyperlink runs ---
Replace-Text "freeCodeCamp – Javascript Algorithms and Data Structures" "freeCodeCamp – Javascript Algorithms and Data Structures"

# --- EDUCATION section: merge "B.Sc" + " " runs (drop spelling/grammar marks) ---
Replace-Text "B.Sc " "B.Sc "

# --- PROFESSIONAL SUMMARY: hyperlink the "CloudFix" mention ---
$rng = $d.Content.Duplicate
$rng.Find.Execute("CloudFix", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Find.Execute("CloudFix", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($rng, "https://cloudfix.com/", $null, $null, "CloudFix") | Out-Null

# --- PROFESSIONAL SUMMARY: $18.3M -> $15M+ ---
Replace-Text "$18.3M " "$15M+ "

# --- WORK EXPERIENCE / Trilogy bullets ---
Replace-Text ", attrition," ", product rollouts,"
Replace-Text "Technical Product Manager. " "Technical Product Manager, Central TPM founding member."
Replace-Text "Led initiatives saving $15M, contributing to CloudFix with $100M in savings." "Led initiatives saving $15M+, contributing to CloudFix with $100M+ in savings."
Replace-Text "Developed tools for spec ordering, stack ranking, and CTO bootcamps." "Developed tools for spec ordering, search, stack ranking, and CTO bootcamps."
Replace-Text "Led 20+ senior SREs in automation and production management." "Coached and managed 20+ senior SREs in automation and production management."

# --- Tissow Technology Ventures header: merge runs ---
Replace-Text "Tissow Technology Ventures, LLP." "Tissow Technology Ventures, LLP."

# --- Humingo bullet: merge runs (no text change) ---
Replace-Text "Led full-cycle architecture and development of Humingo, an e-commerce platform." "Led full-cycle architecture and development of Humingo, an e-commerce platform."

# --- Ticketgoose bullet: merge "Modernized " + "Ticketgoose" ---
Replace-Text "Modernized Ticketgoose" "Modernized Ticketgoose"

# --- Kachyng bullet: payment -> payments, add " and ads" ---
Replace-Text "Steered the architecture and deployment of Kachyng, a PCI-compliant mobile payment platform with single-click checkout." "Steered the architecture and deployment of Kachyng, a PCI-compliant mobile payments platform with single-click checkout and ads."

# --- financial research products bullet ---
Replace-Text "Led a 10-member team in software engineering and SaaS, developing two financial research products." "Led a 10-member team in software engineering and SaaS, developing two software products in financial information research space."

# --- YuMe, Inc. header: merge runs ---
Replace-Text "YuMe, Inc." "YuMe, Inc."

# --- Google Web Toolkit bullet: add missing period ---
Replace-Text "Google Web Toolkit" "Google Web Toolkit."

# --- Footer: merge "For " + "a " + "more detailed resume, please visit " ---
$sec = $d.Sections.First
$ftr = $sec.Footers.Item(1)
$ftr.Range.Find.Execute("For a more detailed resume, please visit ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "For a more detailed resume, please visit ", 2) | Out-Null
